$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1160.6818
$ws.Range("I19").Value = 1062.6471
$ws.Range("K19").Value = 1062.6471
$ws.Range("M19").Value = -887.6470999999999
$ws.Range("H39").Value = 736.2308
$ws.Range("I39").Value = 51.81818
$ws.Range("J39").Value = 4500.5
$ws.Range("K39").Value = 155.45454
$ws.Range("L39").Value = 13501.5
$ws.Range("M39").Value = 140.54546
$ws.Range("N39").Value = -14093.5
$ws.Range("H138").Value = 1159.9
$ws.Range("J138").Value = 760
$ws.Range("L138").Value = 2280
$ws.Range("N138").Value = -12560
$ws.Range("H141").Value = 24999.666
$ws.Range("I141").Value = 24999.666
$ws.Range("K141").Value = 74998.99800000001
$ws.Range("M141").Value = -69818.99800000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2422.1064
$ws.Range("I32").Value = 2422.1064
$ws.Range("K32").Value = 2422.1064
$ws.Range("M32").Value = -2135.1064
$ws.Range("H61").Value = 4857.8335
$ws.Range("I61").Value = 4572.1816
$ws.Range("K61").Value = 4572.1816
$ws.Range("M61").Value = -4360.1816
$ws.Range("H88").Value = 102589.5
$ws.Range("J88").Value = 2859.125
$ws.Range("L88").Value = 2859.125
$ws.Range("N88").Value = -3671.125
$ws.Range("H91").Value = 102589.5
$ws.Range("J91").Value = 2859.125
$ws.Range("L91").Value = 2859.125
$ws.Range("N91").Value = -5667.125
$ws.Range("H110").Value = 4538.9165
$ws.Range("I110").Value = 2809
$ws.Range("J110").Value = 7998.75
$ws.Range("K110").Value = 2809
$ws.Range("L110").Value = 7998.75
$ws.Range("M110").Value = -764
$ws.Range("N110").Value = -12088.75
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H132").Value = 4157.5835
$ws.Range("I132").Value = 4353.909
$ws.Range("K132").Value = 13061.727
$ws.Range("M132").Value = -10531.727
$ws.Range("H136").Value = 4857.8335
$ws.Range("I136").Value = 4572.1816
$ws.Range("K136").Value = 13716.5448
$ws.Range("M136").Value = -11166.5448

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1478.0588
$ws.Range("I20").Value = 1448.2
$ws.Range("J20").Value = 1520.7142
$ws.Range("K20").Value = 1448.2
$ws.Range("L20").Value = 1520.7142
$ws.Range("M20").Value = -1201.2
$ws.Range("N20").Value = -2014.7142
$ws.Range("H107").Value = 2208.2144
$ws.Range("J107").Value = 2998.75
$ws.Range("L107").Value = 2998.75
$ws.Range("N107").Value = -6838.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 425.625
$ws.Range("I5").Value = 321.4
$ws.Range("J5").Value = 599.3333
$ws.Range("K5").Value = 321.4
$ws.Range("L5").Value = 599.3333
$ws.Range("M5").Value = -209.4
$ws.Range("N5").Value = -823.3333
$ws.Range("H8").Value = 1649.3334
$ws.Range("J8").Value = 2999
$ws.Range("L8").Value = 2999
$ws.Range("N8").Value = -3279
$ws.Range("H31").Value = 2372.2
$ws.Range("I31").Value = 2303.6667
$ws.Range("J31").Value = 2475
$ws.Range("K31").Value = 2303.6667
$ws.Range("L31").Value = 2475
$ws.Range("M31").Value = -2008.6667
$ws.Range("N31").Value = -3065
$ws.Range("H34").Value = 2372.2
$ws.Range("I34").Value = 2303.6667
$ws.Range("J34").Value = 2475
$ws.Range("K34").Value = 2303.6667
$ws.Range("L34").Value = 2475
$ws.Range("M34").Value = -2101.6667
$ws.Range("N34").Value = -2879
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 303
$ws.Range("I23").Value = 279.75
$ws.Range("K23").Value = 839.25
$ws.Range("M23").Value = -604.25
$ws.Range("H39").Value = 7221.4443
$ws.Range("J39").Value = 7221.4443
$ws.Range("L39").Value = 21664.3329
$ws.Range("N39").Value = -22252.3329
$ws.Range("H41").Value = 5500
$ws.Range("J41").Value = 5500
$ws.Range("L41").Value = 16500
$ws.Range("N41").Value = -17176
$ws.Range("H52").Value = 6475
$ws.Range("J52").Value = 6475
$ws.Range("L52").Value = 19425
$ws.Range("N52").Value = -19957
$ws.Range("H55").Value = 8551
$ws.Range("I55").Value = 4860
$ws.Range("J55").Value = 10396.5
$ws.Range("K55").Value = 14580
$ws.Range("L55").Value = 31189.5
$ws.Range("M55").Value = -14403
$ws.Range("N55").Value = -31543.5
$ws.Range("H109").Value = 5328.643
$ws.Range("I109").Value = 969
$ws.Range("K109").Value = 2907
$ws.Range("M109").Value = -1867
$ws.Range("H131").Value = 1049
$ws.Range("I131").Value = 1049
$ws.Range("K131").Value = 3147
$ws.Range("M131").Value = 1893

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 594.125
$ws.Range("I97").Value = 536.2857
$ws.Range("K97").Value = 536.2857
$ws.Range("M97").Value = -40.28570000000002
$ws.Range("H102").Value = 2712.5715
$ws.Range("I102").Value = 2613.6155
$ws.Range("J102").Value = 3999
$ws.Range("K102").Value = 2613.6155
$ws.Range("L102").Value = 3999
$ws.Range("M102").Value = -991.6154999999999
$ws.Range("N102").Value = -7243
$ws.Range("H122").Value = 4447.8335
$ws.Range("I122").Value = 4937.6
$ws.Range("J122").Value = 1999
$ws.Range("K122").Value = 14812.8
$ws.Range("L122").Value = 5997
$ws.Range("M122").Value = -12362.8
$ws.Range("N122").Value = -10897
$ws.Range("H126").Value = 3451
$ws.Range("I126").Value = 2425
$ws.Range("K126").Value = 7275
$ws.Range("M126").Value = -4805

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1068.2
$ws.Range("I22").Value = 834
$ws.Range("J22").Value = 1354.4445
$ws.Range("K22").Value = 834
$ws.Range("L22").Value = 1354.4445
$ws.Range("M22").Value = -539
$ws.Range("N22").Value = -1944.4445
$ws.Range("H27").Value = 1068.2
$ws.Range("I27").Value = 834
$ws.Range("J27").Value = 1354.4445
$ws.Range("K27").Value = 834
$ws.Range("L27").Value = 1354.4445
$ws.Range("M27").Value = -727
$ws.Range("N27").Value = -1568.4445
$ws.Range("H40").Value = 3097.5
$ws.Range("I40").Value = 3097.5
$ws.Range("K40").Value = 3097.5
$ws.Range("M40").Value = -2961.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3541.2222
$ws.Range("I122").Value = 3640.9375
$ws.Range("K122").Value = 10922.8125
$ws.Range("M122").Value = -8472.8125
$ws.Range("H132").Value = 4586.591
$ws.Range("I132").Value = 4173.1577
$ws.Range("K132").Value = 12519.4731
$ws.Range("M132").Value = -9989.473099999999
